# "Adds languages to resume"
#
# 1. Merge the " Firebase" + "," runs (removing the stray _GoBack bookmark
#    that sat between them) into a single " Firebase," run.
# 2. Change "React.js, Node.js, and Express.js" to
#    "React.js, Node.js, Express.js, Ggplot2, Plotly, D3.js" (drop the
#    "and", append the data-viz skills that used to live on their own line).
# 3. Rename the "Data Visualization:" label to "Languages:".
# 4. Replace the (now-duplicate) "Ggplot2, Plotly, D3.js" value on the
#    Languages line with the actual language list.

$d = $word.ActiveDocument

# --- Step 1: fold " Firebase" + bookmark + "," into " Firebase," ---
$found1 = $d.Content.Find.Execute(
    "Javascript, Firebase,", $true, $false, $false, $false, $false,
    $true, 1, $false, "Javascript, Firebase,", 2)
Write-Output "Firebase merge executed: $found1"

# --- Step 2: drop "and " before Express.js, tack on the data-viz list ---
$found2 = $d.Content.Find.Execute(
    "React.js, Node.js, and Express.js", $true, $false, $false, $false, $false,
    $true, 1, $false, "React.js, Node.js, Express.js, Ggplot2, Plotly, D3.js", 2)
Write-Output "Programming line updated: $found2"

# --- Step 3: rename "Data Visualization:" label to "Languages:" ---
$found3 = $d.Content.Find.Execute(
    "Data Visualization:", $true, $false, $false, $false, $false,
    $true, 1, $false, "Languages:", 2)
Write-Output "Label renamed: $found3"

# --- Step 4: swap the stale "Ggplot2, Plotly, D3.js" value for languages ---
$found4 = $d.Content.Find.Execute(
    "Languages:" + [char]9 + "Ggplot2, Plotly, D3.js", $true, $false, $false, $false, $false,
    $true, 1, $false, "Languages:" + [char]9 + "English (Fluent), Spanish (Fluent), German (B1 Certified)", 2)
Write-Output "Languages value set: $found4"

Write-Output "Done."
